# Rename placeholders in the "experiment requests" report template sheet:
#   ${experiment.uuid}            -> ${experiment.requestId}
#   ${experiment.experimentStatus} -> ${experiment.requestStatus}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Отчет по заявкам на эксперимент")

$ws.Range("A7").Value = '${experiment.requestId}'
$ws.Range("C7").Value = '${experiment.requestStatus}'
